$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.572.10'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  +0.31%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.696.01'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  +0.43%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '676.17'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -1.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.16'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +0.98%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  +0.77%  '

$ws.Range('E9').Value = '  +1.37%  '

$ws.Range('E10').Value = '  +0.09%  '

$ws.Range('E12').Value = '  +0.83%  '

$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.685.64'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  +0.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.533.58'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  +0.30%  '

$ws.Range('E16').Value = '  +2.32%  '

$ws.Range('E17').Value = '  +1.10%  '

$ws.Range('E18').Value = '  +0.37%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '470.88'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  +0.72%  '

$ws.Range('E20').Value = '  -2.45%  '

$ws.Range('E21').Value = '  +0.99%  '

$ws.Range('E22').Value = '  +1.41%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.842.94'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  +2.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.88'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('E27').Value = '  -0.32%  '

$ws.Range('E28').Value = '  +0.56%  '

$ws.Range('E29').Value = '  +2.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.02'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  +0.58%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.59'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  +0.13%  '

$ws.Range('B32').Value = 'EthereumClassic'

$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.98'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +1.22%  '

$ws.Range('B33').Value = 'Binance-PegBSC-USD'

$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.685.83'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  +0.89%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.162'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('E36').Value = '  +4.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.22'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  +1.75%  '

$ws.Range('E39').Value = '  -1.49%  '

$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('E41').Value = '  +0.23%  '

$ws.Range('B42').Value = 'Mantle'

$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.944'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  +0.37%  '

$ws.Range('B43').Value = 'Monero'

$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '167.07'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  +0.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '46.64'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -2.25%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.76'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  +1.99%  '

$ws.Range('B46').Value = 'InjectiveProtocol'

$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.19'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  +0.88%  '

$ws.Range('B47').Value = 'FLOKI'

$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000279'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  +1.87%  '

$ws.Range('E48').Value = '  +0.04%  '

$ws.Range('E49').Value = '  -2.45%  '

$ws.Range('E50').Value = '  +1.31%  '

$ws.Range('E51').Value = '  +2.06%  '
